# Apply cryptos list price/volume updates (commit: Tue Dec 19 17:39:59 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.986.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.156.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.54%  "
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.09"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.16%  "
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.580"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0904"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.71"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.480.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.14"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.167.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("E17").Value = "  -2.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.858.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.87%  "
$ws.Range("E19").Value = "  -2.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "225.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.93%  "
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.35%  "
$ws.Range("E27").Value = "  +1.04%  "
$ws.Range("E28").Value = "  +2.39%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +10.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "168.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.87"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0792"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.08"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.67%  "
$ws.Range("E35").Value = "  -0.82%  "
$ws.Range("E36").Value = "  +0.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.24"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0325"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.21%  "
$ws.Range("E40").Value = "  -2.52%  "
$ws.Range("E41").Value = "  +3.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "58.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.465"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +16.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.23"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0963"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("E48").Value = "  +8.27%  "
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("E51").Value = "  +0.86%  "
